$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.116.13'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.798.51'
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '222.89'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.32'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('E9').Value = '  +2.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0716'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0925'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').Value = '2.057.70'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '1.789.95'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.71'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.632'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').Value = '34.128.46'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.17'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '246.77'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.24%  '
$ws.Range('D20').Value = '0.0₃0787'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.89'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.86%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '159.29'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  +1.64%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +1.88%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.73'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('D35').Value = '1.414.74'
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.646'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.77%  '
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.944'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +5.11%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '80.21'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.15'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.00%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.94'
$ws.Range('D44').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '107.40'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.03%  '
$ws.Range('D47').Value = '1.956.55'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.96'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').Value = '0.0₆0125'
$ws.Range('E51').Value = '  +3.05%  '
